$wb = $excel.ActiveWorkbook

# --- Sheet "Vowels" ---
$wsVowels = $wb.Worksheets.Item("Vowels")

# Swap the Open/Mid and Front/Cent header labels in row 3
$wsVowels.Range("D3").Value = "Mid"
$wsVowels.Range("F3").Value = "Open"
$wsVowels.Range("I3").Value = "Cent"
$wsVowels.Range("J3").Value = "Front"

# Fill in the counts for row 4 (D4:K4)
$wsVowels.Range("D4").Value = 17
$wsVowels.Range("E4").Value = 1
$wsVowels.Range("F4").Value = 2
$wsVowels.Range("G4").Value = 24
$wsVowels.Range("H4").Value = 20
$wsVowels.Range("I4").Value = 39
$wsVowels.Range("J4").Value = 4
$wsVowels.Range("K4").Value = 21

# --- Sheet "Cons manner" ---
$wsConsManner = $wb.Worksheets.Item("Cons manner")
$wsConsManner.Range("E4").Value = 36

# --- Sheet "Cons place" ---
$wsConsPlace = $wb.Worksheets.Item("Cons place")
$wsConsPlace.Range("E4").ClearContents()
